$wb = $excel.ActiveWorkbook

# --- Sheet "P_valores" ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.2899281983818995
$wsP.Range("D2").Value = 0.3666530656187619
$wsP.Range("E2").Value = 0.8129955704105249
$wsP.Range("F2").Value = 0.8030672424573058

$wsP.Range("B3").Value = 0.2899281983818995
$wsP.Range("D3").Value = 0.7970147502260116
$wsP.Range("E3").Value = 0.2486091604779133
$wsP.Range("F3").Value = 0.3488216826859514

$wsP.Range("B4").Value = 0.3666530656187619
$wsP.Range("C4").Value = 0.7970147502260116
$wsP.Range("E4").Value = 0.2654646117753956
$wsP.Range("F4").Value = 0.1588276240389643

$wsP.Range("B5").Value = 0.8129955704105249
$wsP.Range("C5").Value = 0.2486091604779133
$wsP.Range("D5").Value = 0.2654646117753956
$wsP.Range("F5").Value = 0.6754291076950425

$wsP.Range("B6").Value = 0.8030672424573058
$wsP.Range("C6").Value = 0.3488216826859514
$wsP.Range("D6").Value = 0.1588276240389643
$wsP.Range("E6").Value = 0.6754291076950425

# --- Sheet "Estadisticos_DM" ---
$wsE = $wb.Worksheets.Item("Estadisticos_DM")

$wsE.Range("C2").Value = 1.099903333524128
$wsE.Range("D2").Value = 0.9329590579627471
$wsE.Range("E2").Value = -0.2410700623706254
$wsE.Range("F2").Value = 0.25415951958384

$wsE.Range("B3").Value = -1.099903333524128
$wsE.Range("D3").Value = 0.2621617886823465
$wsE.Range("E3").Value = -1.203847019887192
$wsE.Range("F3").Value = -0.9693468040957722

$wsE.Range("B4").Value = -0.9329590579627471
$wsE.Range("C4").Value = -0.2621617886823465
$wsE.Range("E4").Value = -1.159961669520672
$wsE.Range("F4").Value = -1.488379851961763

$wsE.Range("B5").Value = 0.2410700623706254
$wsE.Range("C5").Value = 1.203847019887192
$wsE.Range("D5").Value = 1.159961669520672
$wsE.Range("F5").Value = 0.4276243293568978

$wsE.Range("B6").Value = -0.25415951958384
$wsE.Range("C6").Value = 0.9693468040957722
$wsE.Range("D6").Value = 1.488379851961763
$wsE.Range("E6").Value = -0.4276243293568978

$wb.Save()
